{"js": "// Office.js (Word JavaScript API) script.\n// The visible-text effect of the diff is: the phrase \"puzzle-board\" (as in\n// \"... of the puzzle-board elements in the resetPuzzlePieces function ...\")\n// is changed to \"drop zone\" (so it reads \"... of the drop zone elements in\n// the resetPuzzlePieces function ...\"). The cursor/_GoBack bookmark simply\n// follows the last text edit, which Word manages automatically.\n\nconst results = context.document.body.search(\"puzzle-board\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"drop zone\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# The visible-text effect of the diff is: \"puzzle-board\" is changed to\n# \"drop zone\" in the sentence \"... of the puzzle-board elements in the\n# resetPuzzlePieces function ...\" -> \"... of the drop zone elements in the\n# resetPuzzlePieces function ...\". The _GoBack bookmark simply tracks the\n# most recent edit location, which Word manages automatically.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"puzzle-board\"\n$find.Replacement.Text = \"drop zone\"\n$find.Forward = $true\n$find.Wrap = 1        # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n# wdReplaceAll = 2\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n"}
